$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Consumer demand assignment correction -------------------------------
# The agent rows (agent_id/role/working_capital/selling_price) need to be
# rotated: row2 <- old row3, row3 <- old row4, row4 <- old row2, with the
# wrapped-around row also getting a freshly assigned working_capital and a
# numeric (not text) selling_price.

# Stash old row 2 (A2:K2) in a scratch row so it is not lost.
$ws.Range("A2:K2").Copy()
$ws.Range("A31").PasteSpecial(-4163)

# row2 <- old row3
$ws.Range("A3:K3").Copy()
$ws.Range("A2").PasteSpecial(-4163)

# row3 <- old row4
$ws.Range("A4:K4").Copy()
$ws.Range("A3").PasteSpecial(-4163)

# row4 <- old row2 (from scratch)
$ws.Range("A31:K31").Copy()
$ws.Range("A4").PasteSpecial(-4163)

# Corrected demand assignment for the wrapped-around agent.
$ws.Range("C4").Value = 150
$ws.Range("D4").Value = 5.5

# Clean up the scratch row.
$ws.Range("A31:K31").Clear()

# --- Remove the now-superfluous trailing blank row -----------------------
$ws.Rows.Item(29).Delete()

# --- Restore the selection Excel leaves after these edits ----------------
$ws.Range("K3:K4").Select() | Out-Null
